$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet (SPY -> EWZ) -------------------------------------------------
$ws.Name = "20161014 EWZ IV smile data"

# --- Update the defined name to match new sheet name + expanded range ---------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_20160928_UNG") {
        $n.RefersTo = "='20161014 EWZ IV smile data'!`$A`$1:`$D`$53"
    }
}

# --- Insert rows to make room for the new fields -------------------------------
# New "ticker" row right after option_type: duplicate row 5 (which already
# carries the styled-but-empty column-F formatting) and insert the copy at
# row 6, pushing the original expiration_yy/mm/dd rows down by one while
# keeping the F-column styling contiguous from row 5.
$ws.Rows(5).Copy()
$ws.Rows(6).Insert()

# The new buy/write fields extend the table: everything through row 17
# already lines up 1:1 with existing rows, so just append 6 fresh rows at
# the bottom (18-23) to hold the remaining entry_date_*_write fields and
# the relocated historical_volatility row.
for ($i = 0; $i -lt 6; $i++) {
    $ws.Rows(18).Insert()
}

# --- Write the final field/value layout ----------------------------------------
$ws.Range("A1").Value = "FIELD"
$ws.Range("B1").Value = "VALUE"

$ws.Range("A2").Value = "smile_file"
$ws.Range("B2").Value = "20161014 EWZ"

$ws.Range("A3").Value = "smile_file_path"
$ws.Range("B3").Value = "E:\\Datos\\bolsa\\cuenta personal\\analisis de valores\\Trades activos\\Scanning\\20161014"

$ws.Range("A4").Value = "option_type"
$ws.Range("B4").Value = "put"

$ws.Range("A5").Value = "ticker"
$ws.Range("B5").Value = "EWZ"

$ws.Range("A6").Value = "expiration_yy"
$ws.Range("B6").Value = 2016

$ws.Range("A7").Value = "expiration_mm"
$ws.Range("B7").Value = 10

$ws.Range("A8").Value = "expiration_dd"
$ws.Range("B8").Value = 21

$ws.Range("A9").Value = "strike_buy"
$ws.Range("B9").Value = 35

$ws.Range("A10").Value = "entry_date_yy_buy"
$ws.Range("B10").Value = 2016

$ws.Range("A11").Value = "entry_date_mm_buy"
$ws.Range("B11").Value = 10

$ws.Range("A12").Value = "entry_date_dd_buy"
$ws.Range("B12").Value = 14

$ws.Range("A13").Value = "entry_date_hh_buy"
$ws.Range("B13").Value = 9

$ws.Range("A14").Value = "entry_date_MM_buy"
$ws.Range("B14").Value = 50

$ws.Range("A15").Value = "entry_date_ss_buy"
$ws.Range("B15").Value = 23

$ws.Range("A16").Value = "strike_write"
$ws.Range("B16").Value = 36

$ws.Range("A17").Value = "entry_date_yy_write"
$ws.Range("B17").Value = 2016

$ws.Range("A18").Value = "entry_date_mm_write"
$ws.Range("B18").Value = 10

$ws.Range("A19").Value = "entry_date_dd_write"
$ws.Range("B19").Value = 14

$ws.Range("A20").Value = "entry_date_hh_write"
$ws.Range("B20").Value = 9

$ws.Range("A21").Value = "entry_date_MM_write"
$ws.Range("B21").Value = 58

$ws.Range("A22").Value = "entry_date_ss_write"
$ws.Range("B22").Value = 14

$ws.Range("A23").Value = "historical_volatility"
$ws.Range("B23").Value = 0

# --- Column F only keeps its styled-but-empty cells through row 17 -------------
$ws.Range("F18:F23").Clear()

# --- Selection marker matches the new last-used cell ---------------------------
$ws.Range("B23").Select()
